# Weekly fruit/vegetable price update.
# A new daily price record for "Acelga" (Vega Modelo de Temuco) was
# collected and inserted into the existing date-ordered table at row 166,
# pushing the previous rows 166-182 down to 167-183.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 166 (shifts old rows 166-182 down to 167-183,
# carrying their values/formatting with them).
$ws.Rows.Item(166).Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A166").Value = 10
$ws.Range("B166").Value = "Vega Modelo de Temuco"
$ws.Range("C166").Value = "La Araucanía"
$ws.Range("D166").Value = 44449
$ws.Range("E166").Value = 9
$ws.Range("F166").Value = 100112009
$ws.Range("G166").Value = "Acelga"
$ws.Range("H166").Value = "Sin especificar"
$ws.Range("I166").Value = "Primera"
$ws.Range("J166").Value = 65
$ws.Range("K166").Value = 8000
$ws.Range("L166").Value = 8000
$ws.Range("M166").Value = 8000
$ws.Range("N166").Value = "$/docena de atados (12 kilos)"
$ws.Range("O166").Value = "Provincia de Cautín"
$ws.Range("P166").Value = 667
$ws.Range("Q166").Value = 12
$ws.Range("R166").Value = "Hortaliza"
